# Append new trade row (Trade #14) to the "All Trades" and "base_strategy"
# sheets, mirroring a new open-trade log entry at 2026-02-16 22:53:24.

$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "base_strategy")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # New row goes right after the current last row (row 14 -> row 15).
    $row = 15

    # Date column: force text so "2026-02-16" isn't auto-converted to a
    # date serial by Excel's literal-entry parsing, then clear the
    # number-format style back off so no style index lingers on the cell.
    $dateCell = $ws.Cells.Item($row, 2)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = "2026-02-16"
    $dateCell.Style = "Normal"

    $ws.Cells.Item($row, 1).Value = 14                 # A: Trade #
    $ws.Cells.Item($row, 3).Value = "22:53:24"          # C: Time
    $ws.Cells.Item($row, 4).Value = "base_strategy"     # D: Strategy
    $ws.Cells.Item($row, 5).Value = "DOWN"              # E: Side
    $ws.Cells.Item($row, 6).Value = 49.999998           # F: Entry Price
    $ws.Cells.Item($row, 7).Value = ""                  # G: Exit Price
    $ws.Cells.Item($row, 8).Value = "OPEN"              # H: Status
    $ws.Cells.Item($row, 9).Value = 0                   # I: P&L %
    $ws.Cells.Item($row, 10).Value = 0                  # J: P&L $
    $ws.Cells.Item($row, 11).Value = 100                # K: Capital After
    $ws.Cells.Item($row, 12).Value = 0                  # L: Entry Slippage (bps)
    $ws.Cells.Item($row, 13).Value = 0                  # M: Exit Slippage (bps)
    $ws.Cells.Item($row, 14).Value = 0.6                # N: Confidence
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps" # O: Entry Reason
    $ws.Cells.Item($row, 16).Value = ""                 # P: Exit Reason
    $ws.Cells.Item($row, 17).Value = 0                  # Q: Duration (min)
}
